# Update the 15 lattice-multiplication problems in the single 5x3 table.
# Each cell holds a run with 5 <w:t> pieces separated by <w:br/> line breaks:
#   "AB x CD"
#   "  C    D"   (digits of the second factor, spaced)
#   "  ----"
#   "A|    |"    (first digit of the first factor)
#   "B|    |"    (second digit of the first factor)
#
# We rewrite each table cell's Range.Text wholesale (using a vertical-tab
# character to represent the <w:br/> line breaks that Word uses inside a
# cell), which keeps the original run/paragraph formatting (sz=32) intact
# and avoids any ambiguity from overlapping Find/Replace matches such as
# "7|    |" -> "3|    |" followed by "3|    |" -> "1|    |".

$d = $word.ActiveDocument
$tbl = $d.Tables(1)
$vt = [char]11

$tbl.Cell(1,1).Range.Text = "31 x 35" + $vt + "  3    5" + $vt + "  ----" + $vt + "3|    |" + $vt + "1|    |"
$tbl.Cell(1,2).Range.Text = "83 x 24" + $vt + "  2    4" + $vt + "  ----" + $vt + "8|    |" + $vt + "3|    |"
$tbl.Cell(1,3).Range.Text = "70 x 81" + $vt + "  8    1" + $vt + "  ----" + $vt + "7|    |" + $vt + "0|    |"

$tbl.Cell(2,1).Range.Text = "71 x 33" + $vt + "  3    3" + $vt + "  ----" + $vt + "7|    |" + $vt + "1|    |"
$tbl.Cell(2,2).Range.Text = "13 x 52" + $vt + "  5    2" + $vt + "  ----" + $vt + "1|    |" + $vt + "3|    |"
$tbl.Cell(2,3).Range.Text = "23 x 97" + $vt + "  9    7" + $vt + "  ----" + $vt + "2|    |" + $vt + "3|    |"

$tbl.Cell(3,1).Range.Text = "71 x 13" + $vt + "  1    3" + $vt + "  ----" + $vt + "7|    |" + $vt + "1|    |"
$tbl.Cell(3,2).Range.Text = "93 x 93" + $vt + "  9    3" + $vt + "  ----" + $vt + "9|    |" + $vt + "3|    |"
$tbl.Cell(3,3).Range.Text = "42 x 76" + $vt + "  7    6" + $vt + "  ----" + $vt + "4|    |" + $vt + "2|    |"

$tbl.Cell(4,1).Range.Text = "76 x 47" + $vt + "  4    7" + $vt + "  ----" + $vt + "7|    |" + $vt + "6|    |"
$tbl.Cell(4,2).Range.Text = "34 x 37" + $vt + "  3    7" + $vt + "  ----" + $vt + "3|    |" + $vt + "4|    |"
$tbl.Cell(4,3).Range.Text = "69 x 65" + $vt + "  6    5" + $vt + "  ----" + $vt + "6|    |" + $vt + "9|    |"

$tbl.Cell(5,1).Range.Text = "14 x 65" + $vt + "  6    5" + $vt + "  ----" + $vt + "1|    |" + $vt + "4|    |"
$tbl.Cell(5,2).Range.Text = "20 x 76" + $vt + "  7    6" + $vt + "  ----" + $vt + "2|    |" + $vt + "0|    |"
$tbl.Cell(5,3).Range.Text = "86 x 56" + $vt + "  5    6" + $vt + "  ----" + $vt + "8|    |" + $vt + "6|    |"
